# Update "想去人数" (wanted-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - rows 3-6, column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 206
$wsExhibit.Range("F4").Value = 2441
$wsExhibit.Range("F5").Value = 33
$wsExhibit.Range("F6").Value = 530

# Sheet "全部类型" (All types) - rows 5-8, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 206
$wsAll.Range("F6").Value = 2441
$wsAll.Range("F7").Value = 33
$wsAll.Range("F8").Value = 530
